# Update microstates file to v1.3.2.
# Resonance structures and most geometric isomers eliminated:
#   SM19_micro003, SM19_micro009, SM19_micro011 and SM19_micro017 (with
#   their 2D-depiction pictures) are removed, and the title cell is
#   updated from the old CSV filename to "Microstate List".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Update the title cell (A1) from the old CSV filename to the new title.
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "Microstate List"

# ---------------------------------------------------------------------------
# 2) Remove the 2D-depiction pictures belonging to the four eliminated
#    microstates first (before touching positions/rows of anything else).
#    Picture "N" is the depiction that originally sat next to worksheet
#    row (N + 2): Picture 3 -> row 5 (SM19_micro003), Picture 9 -> row 11
#    (SM19_micro009), Picture 11 -> row 13 (SM19_micro011), Picture 17 ->
#    row 19 (SM19_micro017).
# ---------------------------------------------------------------------------
$picturesToRemove = @("Picture 3", "Picture 9", "Picture 11", "Picture 17")
foreach ($picName in $picturesToRemove) {
    $ws.Shapes.Item($picName).Delete()
}

# ---------------------------------------------------------------------------
# 3) Re-anchor every surviving picture to the worksheet row it will occupy
#    once the now-empty rows are deleted, while all of the original rows
#    still exist (so row tops/heights used for the target position are
#    still valid).
# ---------------------------------------------------------------------------
$survivingPictureNewRow = [ordered]@{
    "Picture 1"  = 3
    "Picture 2"  = 4
    "Picture 4"  = 5
    "Picture 5"  = 6
    "Picture 6"  = 7
    "Picture 7"  = 8
    "Picture 8"  = 9
    "Picture 10" = 10
    "Picture 12" = 11
    "Picture 13" = 12
    "Picture 14" = 13
    "Picture 15" = 14
    "Picture 16" = 15
    "Picture 18" = 16
    "Picture 19" = 17
}

foreach ($picName in $survivingPictureNewRow.Keys) {
    $targetRow = $survivingPictureNewRow[$picName]
    $shp = $ws.Shapes.Item($picName)
    $shp.Top = $ws.Rows($targetRow.ToString() + ":" + $targetRow.ToString()).Top
    $shp.Left = 0
}

# ---------------------------------------------------------------------------
# 4) Remove the worksheet rows for the four eliminated microstates
#    (SM19_micro003 -> row 5, SM19_micro009 -> row 11, SM19_micro011 -> row
#    13, SM19_micro017 -> row 19). Delete from the bottom up so earlier row
#    numbers stay valid while iterating.
# ---------------------------------------------------------------------------
$rowsToRemove = @(19, 13, 11, 5)
foreach ($r in $rowsToRemove) {
    $ws.Rows($r.ToString() + ":" + $r.ToString()).Delete()
}
